# New features for updated outcome measures: add three new settings rows
# ("desat severe threshold", "expected_sampling_rate (sec)" and
# "artifact duration threshold (sec)") to the "settings" sheet, and bump
# "desat spike" from -5 to -10.
#
# Work from the bottom of the sheet upwards so inserted rows don't shift
# the row numbers we still need to touch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Insert "expected_sampling_rate (sec)" = 4 before the current row 9
# ("night duration bin size (hours)"), i.e. right after "minimum night
# duration (hours)".
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "expected_sampling_rate (sec)"
$ws.Cells.Item(9, 2).Value = 4

# Insert "desat severe threshold" = 85 before the current row 4
# ("desat spike"), i.e. right after "desat threshold".
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = "desat severe threshold"
$ws.Cells.Item(4, 2).Value = 85

# "desat spike" (now shifted down to row 5) changes from -5 to -10.
$ws.Cells.Item(5, 2).Value = -10

# Append "artifact duration threshold (sec)" = 30 as the new last row.
$ws.Cells.Item(15, 1).Value = "artifact duration threshold (sec)"
$ws.Cells.Item(15, 2).Value = 30

$ws.Range("A16").Select()
